# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.17 = 12161.74 pesos`n✅ 12161.74 pesos = 3.15 = 960.47 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 3837.03
$ws2.Range("N12").Value = 3862
$ws2.Range("O12").Value = 305
